$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F and G columns for rows 334-368 (AgTests / AgPosit revisions)
$ws.Range("F334").Value = 194969
$ws.Range("G334").Value = 3427
$ws.Range("F335").Value = 130563
$ws.Range("G335").Value = 2986
$ws.Range("F336").Value = 101915
$ws.Range("G336").Value = 3358
$ws.Range("F337").Value = 104151
$ws.Range("G337").Value = 2980
$ws.Range("F338").Value = 226630
$ws.Range("G338").Value = 3191
$ws.Range("F340").Value = 380718
$ws.Range("G340").Value = 3262
$ws.Range("F341").Value = 295365
$ws.Range("G341").Value = 3672
$ws.Range("F342").Value = 179413
$ws.Range("G342").Value = 3077
$ws.Range("F343").Value = 133537
$ws.Range("G343").Value = 2980
$ws.Range("F344").Value = 135985
$ws.Range("G344").Value = 2489
$ws.Range("F345").Value = 290482
$ws.Range("F346").Value = 666492
$ws.Range("F347").Value = 340341
$ws.Range("G347").Value = 2887
$ws.Range("F348").Value = 231837
$ws.Range("G348").Value = 3238
$ws.Range("F349").Value = 159802
$ws.Range("G349").Value = 2750
$ws.Range("F350").Value = 127604
$ws.Range("G350").Value = 2974
$ws.Range("F351").Value = 150116
$ws.Range("G351").Value = 2821
$ws.Range("F352").Value = 306342
$ws.Range("G352").Value = 3546
$ws.Range("F353").Value = 717568
$ws.Range("G353").Value = 5247
$ws.Range("F354").Value = 304218
$ws.Range("G354").Value = 2779
$ws.Range("F355").Value = 221911
$ws.Range("G355").Value = 3443
$ws.Range("F356").Value = 160337
$ws.Range("G356").Value = 2895
$ws.Range("F357").Value = 138300
$ws.Range("G357").Value = 3018
$ws.Range("F358").Value = 157658
$ws.Range("G358").Value = 2600
$ws.Range("F359").Value = 319936
$ws.Range("G359").Value = 3345
$ws.Range("F360").Value = 738582
$ws.Range("G360").Value = 5035
$ws.Range("F361").Value = 329258
$ws.Range("G361").Value = 2580
$ws.Range("F362").Value = 223568
$ws.Range("G362").Value = 3077
$ws.Range("F363").Value = 184926
$ws.Range("G363").Value = 2713
$ws.Range("F364").Value = 163439
$ws.Range("G364").Value = 2401
$ws.Range("F365").Value = 177342
$ws.Range("G365").Value = 2354
$ws.Range("F366").Value = 328448
$ws.Range("G366").Value = 2773
$ws.Range("F367").Value = 734786
$ws.Range("G367").Value = 3723
$ws.Range("F368").Value = 336673
$ws.Range("G368").Value = 2227

# Row 369: replace B/C/D/F/G with corrected values (E unchanged)
$ws.Range("B369").Value = 325993
$ws.Range("C369").Value = 11222
$ws.Range("D369").Value = 2207
$ws.Range("F369").Value = 221336
$ws.Range("G369").Value = 2402

# New row 370: data for 2021-03-09 (serial 44264)
$ws.Range("A370").Value = 44264
$ws.Range("B370").Value = 329593
$ws.Range("C370").Value = 16361
$ws.Range("D370").Value = 3600
$ws.Range("E370").Value = 8146
$ws.Range("F370").Value = 154883
$ws.Range("G370").Value = 1681
